$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Row 1 header: shift B->A, C->B (new C header set later, after SymbolDb lookups) ---
$ws.Range("A1").Value = $ws.Range("B1").Value()
$ws.Range("B1").Value = $ws.Range("C1").Value()
$ws.Range("B1").Copy()
$ws.Range("A1").PasteSpecial(-4122)
$excel.CutCopyMode = $false

# --- Rows 2-44: shift B->A (coin ticker), C->B (symbol_id), clear old numbering style on A ---
for ($r = 2; $r -le 44; $r++) {
  $coin = $ws.Cells.Item($r, 2).Value()
  $symid = $ws.Cells.Item($r, 3).Value()
  $ws.Cells.Item($r, 1).Value = $coin
  $ws.Cells.Item($r, 1).Style = "Normal"
  $ws.Cells.Item($r, 2).Value = $symid
}

# --- Column C: clear stale symbol_id carried over, for rows with no SymbolDb mapping ---
$ws.Cells.Item(2, 3).ClearContents()
$ws.Cells.Item(9, 3).ClearContents()
$ws.Cells.Item(19, 3).ClearContents()
$ws.Cells.Item(25, 3).ClearContents()
$ws.Cells.Item(26, 3).ClearContents()
$ws.Cells.Item(27, 3).ClearContents()
$ws.Cells.Item(34, 3).ClearContents()
$ws.Cells.Item(39, 3).ClearContents()
$ws.Cells.Item(42, 3).ClearContents()
$ws.Cells.Item(44, 3).ClearContents()

# --- Column C: new "SymbolDb" lookup values, in original authoring order ---
$ws.Cells.Item(4, 3).Value = " ALICEBUSD"
$ws.Cells.Item(32, 3).Value = " HNTUSDT"
$ws.Cells.Item(43, 3).Value = " RNDRUSDT"
$ws.Cells.Item(10, 3).Value = " EQZUSDT"
$ws.Cells.Item(40, 3).Value = " FORMUSDT"
$ws.Cells.Item(33, 3).Value = " HYVEUSDT"
$ws.Cells.Item(3, 3).Value = " SUPERUSDT"
$ws.Cells.Item(37, 3).Value = " ACHUSD"
$ws.Cells.Item(12, 3).Value = " MLNUSDT"
$ws.Cells.Item(41, 3).Value = " CGGUSDT"
$ws.Cells.Item(18, 3).Value = " ADAUSDT"
$ws.Cells.Item(21, 3).Value = " HBARUSDT"
$ws.Cells.Item(14, 3).Value = " RAYBUSD"
$ws.Cells.Item(23, 3).Value = " AXSUSDT"
$ws.Cells.Item(36, 3).Value = " SUSHIUSDT"
$ws.Cells.Item(28, 3).Value = " DUSKUSDT"
$ws.Cells.Item(5, 3).Value = " CROUSDT"
$ws.Cells.Item(35, 3).Value = " SANDUSDT"
$ws.Cells.Item(20, 3).Value = " LTCUSDT"
$ws.Cells.Item(6, 3).Value = " KCSUSDT"
$ws.Cells.Item(7, 3).Value = " POLSUSDT"
$ws.Cells.Item(17, 3).Value = " DATABUSD"
$ws.Cells.Item(29, 3).Value = " FARMUSDT"
$ws.Cells.Item(15, 3).Value = " MANAUSDT"
$ws.Cells.Item(30, 3).Value = " GENSUSDT"
$ws.Cells.Item(11, 3).Value = " STEPUSD"
$ws.Cells.Item(31, 3).Value = " VRAUSDT"
$ws.Cells.Item(24, 3).Value = " ONEUSDT"
$ws.Cells.Item(22, 3).Value = " PHAUSDT"
$ws.Cells.Item(38, 3).Value = " SHIBUSDT"
$ws.Cells.Item(8, 3).Value = " FTMUSDT"
$ws.Cells.Item(16, 3).Value = " BNBUSDT"
$ws.Cells.Item(13, 3).Value = " LOCGUSDT"

# --- Row 1: new column C header (added last so it gets the final shared-string slot) ---
$ws.Range("C1").Value = "SymbolDb"
$ws.Range("C1").Style = "Normal"

# --- Column widths ---
$ws.Columns.Item(1).ColumnWidth = 9.42578125
$ws.Columns.Item(2).ColumnWidth = 29.42578125

# --- Selection ---
$ws.Range("A1:A1048576").Select()
